# Update odds values in row 2 of the active worksheet to reflect the
# latest FlashScore data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 3
$ws.Range("J2").Value = 3
$ws.Range("K2").Value = 2.2
$ws.Range("L2").Value = 3.4

$ws.Range("Q2").Value = 1.83
$ws.Range("R2").Value = 1.98
$ws.Range("S2").Value = 1.36
$ws.Range("T2").Value = 3

$ws.Range("Z2").Value = 23
$ws.Range("AA2").Value = 19
$ws.Range("AB2").Value = 26

$ws.Range("AH2").Value = 11
$ws.Range("AI2").Value = 15
$ws.Range("AK2").Value = 29
$ws.Range("AL2").Value = 21

$ws.Range("AO2").Value = 13

$ws.Range("AT2").Value = 3
